$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new IFERROR formulas / array formula block (rows 14-16, cols K/M/N)
$ws.Range("K14").Formula = "=IFERROR(K12,3)"
$ws.Range("K15").Formula = "=IFERROR(K14,1)"
$ws.Range("M14:N16").FormulaArray = "=IFERROR(I16:J18,1)"

# Update the saved selection to L16:L17 (active cell L16)
$ws.Range("L16:L17").Select()
